# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (strikeout) values per row, recomputed from the underlying
# game logs (replacing the old "Strike#" derived values).
$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 2
    11 = 3
    12 = 1
    13 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 1
    26 = 0
    27 = 2
    28 = 1
    29 = 0
    30 = 1
    31 = 0
    32 = 0
    33 = 1
    34 = 3
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 1
    40 = 1
    41 = 1
    42 = 2
    45 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
